# Rename the header "Genre" (column D) to "Gendre" and move the active
# selection to D1, matching the author's "update the function of leaderboard"
# commit (the underlying gender column values are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Gendre"

$ws.Range("D1").Select()
